$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date style already used in column A (e.g. A17) onto the new date cells
$ws.Range("A17").Copy()
$ws.Range("A18:A19").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row 18
$ws.Cells.Item(18, 1).Value = 43839
$ws.Cells.Item(18, 2).Value = "Design the navigaction"
$ws.Cells.Item(18, 4).Value = "design the navigation"

# Row 19
$ws.Cells.Item(19, 1).Value = 43840
$ws.Cells.Item(19, 2).Value = "design for the test concept"

$ws.Range("D19").Select()
